# Problem 10 draft update:
#  - append a brand-new "bbbba" block (rows 73-77) with the DP table values filled in,
#    copied down from the existing "bbbba" block (rows 66-70)
#  - simplify the original block (rows 66-70) into a smaller "ba" test case

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Rows 73-77: new block appended below row 70, built first while rows 66-70
# still hold their original "bbbba" formatting/styles to copy from.
# ---------------------------------------------------------------------------

# Header / title row (merged like the other section headers above).
$ws.Range("B66:AA66").Copy()
$ws.Range("B73:AA73").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B73:AA73").Merge()
$ws.Range("B73").Value = "bbbba"

# Row 74 spells "bbbba" across B74:F74 (style 1, same as row 67's letters).
$ws.Range("B67").Copy()
$ws.Range("B74:E74").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E67").Copy()
$ws.Range("F74").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B74").Value = "b"
$ws.Range("C74").Value = "b"
$ws.Range("D74").Value = "b"
$ws.Range("E74").Value = "b"
$ws.Range("F74").Value = "a"

# Row 75: label ".*" (style 1) then the DP row of values 1,2,3,4 (style 3).
$ws.Range("A68").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B68").Copy()
$ws.Range("B75:E75").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A75").Value = ".*"
$ws.Range("B75").Value = 1
$ws.Range("C75").Value = 2
$ws.Range("D75").Value = 3
$ws.Range("E75").Value = 4

# Row 76: label "a*" (style 1) then values 1,2,3,4 (style 3) starting at C76.
$ws.Range("A69").Copy()
$ws.Range("A76").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B68").Copy()
$ws.Range("C76:F76").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A76").Value = "a*"
$ws.Range("C76").Value = 1
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 3
$ws.Range("F76").Value = 4

# Row 77: label "a" (style 1) then the final values 1,2,3 (style 7, like row 56/57).
$ws.Range("A70").Copy()
$ws.Range("A77").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C56").Copy()
$ws.Range("C77:E77").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A77").Value = "a"
$ws.Range("C77").Value = 1
$ws.Range("D77").Value = 2
$ws.Range("E77").Value = 3

# ---------------------------------------------------------------------------
# Rows 66-70: shrink the original "bbbba" block down to a "ba" block.
# ---------------------------------------------------------------------------

# Header cell text: "bbbba" -> "ba" (merged title row B66:AA66 keeps its style).
$ws.Range("B66").Value = "ba"

# Row 67 used to spell "bbb"+"a" across B67:E67; it now spells "b"+"a".
$ws.Range("C67").Value = "a"
$ws.Range("D67").Clear()
$ws.Range("E67").Clear()

# Row 68 had three empty data cells (B68:D68); only B68 remains.
$ws.Range("C68").Clear()
$ws.Range("D68").Clear()

# Row 69 gains an empty styled data cell at C69 (style 3, copied from B68).
$ws.Range("B68").Copy()
$ws.Range("C69").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 70's lone empty data cell moves from E70 to C70 (style 3).
$ws.Range("E70").Clear()
$ws.Range("B68").Copy()
$ws.Range("C70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author left off.
$ws.Range("B75").Select()
